# TaskList.xlsx refactor: clear the per-row "Result"/metric figures for the
# UAD/CPP-P-026A (row 13) and UL/570-P-2901 (row 15) task rows, leaving the
# cell formatting (styles) intact but removing their values — matching the
# author's "Major Refactoring all the column names with seperate
# definitions class" commit, which wiped the stale computed figures for
# these two rows while keeping everything else as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("task_list_1")

# Row 13 (Tag CPP-P-026A): clear Result, IT_energy, IT_ghg_reduction,
# IT_ghg_reduction_percent, IT_ghg_cost, VSD_energy, VSD_ghg_reduction,
# VSD_ghg_reduction_percent and VSD_ghg_cost — but keep Perform (D13),
# IT_Annual_Spendings (I13), VSD_Annual_Spendings (N13) and Remarks (P13)
# untouched (they were already blank).
$ws.Range("E13:H13").ClearContents()
$ws.Range("J13:M13").ClearContents()
$ws.Range("O13").ClearContents()

# Row 15 (Tag 570-P-2901): same treatment.
$ws.Range("E15:H15").ClearContents()
$ws.Range("J15:M15").ClearContents()
$ws.Range("O15").ClearContents()

# The author's cursor/selection ended up on D15 when they saved.
$ws.Range("D15").Select()
